$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a new worksheet named "2022-Q4" positioned before "2022-Q3".
#    Duplicate "总计" (rather than Worksheets.Add) so the new sheet
#    inherits the same sheetPr/pageMargins/sheetFormatPr defaults, then
#    wipe its content.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")
$wsQ3 = $wb.Worksheets.Item("2022-Q3")
$wsTotal.Copy($null, $wsTotal)
$wsQ4 = $wb.Worksheets.Item(2)
$wsQ4.Name = "2022-Q4"
$wsQ4.Cells.Clear()

# ---------------------------------------------------------------------
# 2. Populate the new "2022-Q4" sheet with fund holdings data
# ---------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $wsQ4.Cells.Item(1, 2 + $i).Value = $headers[$i]
}

$data = @(
    @("016174", "汇丰晋信策略优选混合A", "1.78", "74.92", "2.50", "0.0445", 4),
    @("519615", "银河君尚灵活配置混合I", "1.83", "38.98", "1.10", "0.0201", 1),
    @("519613", "银河君尚灵活配置混合A", "1.17", "38.98", "1.10", "0.0129", 1),
    @("016175", "汇丰晋信策略优选混合C", "0.40", "74.92", "2.50", "0.0100", 4),
    @("970073", "东证融汇成长优选混合A", "0.38", "89.59", "0.80", "0.0030", 9),
    @("519614", "银河君尚灵活配置混合C", "0.16", "38.98", "1.10", "0.0018", 1),
    @("015921", "申万菱信国证2000指数增强A", "0.21", "94.00", "0.52", "0.0011", 5),
    @("970074", "东证融汇成长优选混合C", "0.11", "89.59", "0.80", "0.0009", 9),
    @("015922", "申万菱信国证2000指数增强C", "0.08", "94.00", "0.52", "0.0004", 5)
)

$lastRow = 1 + $data.Length

# Columns B (基金代码) and D:G (基金规模/股票总仓位/仓位占比/持有市值) must stay
# text (preserve leading zeros / trailing zeros exactly as scraped), so force
# them to Text format before writing, then strip the number-format again so
# the cells end up with no explicit style (matching plain scraped cells).
$wsQ4.Range("B2:B$lastRow").NumberFormat = "@"
$wsQ4.Range("D2:G$lastRow").NumberFormat = "@"

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = 2 + $r
    $rowData = $data[$r]
    $wsQ4.Cells.Item($row, 1).Value = $r
    $wsQ4.Cells.Item($row, 2).Value = $rowData[0]
    $wsQ4.Cells.Item($row, 3).Value = $rowData[1]
    $wsQ4.Cells.Item($row, 4).Value = $rowData[2]
    $wsQ4.Cells.Item($row, 5).Value = $rowData[3]
    $wsQ4.Cells.Item($row, 6).Value = $rowData[4]
    $wsQ4.Cells.Item($row, 7).Value = $rowData[5]
    $wsQ4.Cells.Item($row, 8).Value = $rowData[6]
}

$wsQ4.Range("B2:B$lastRow").ClearFormats()
$wsQ4.Range("D2:G$lastRow").ClearFormats()

# Apply the same styling used by the "总计" sheet header/index cells
# (bold font + border, matching style used in B1:D1 / A2 of 总计)
$wsTotal.Range("B1:D1").Copy()
$wsQ4.Range("B1:H1").PasteSpecial(-4122)
$wsTotal.Range("A2").Copy()
$wsQ4.Range("A2:A$lastRow").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. Update the "总计" sheet: push the existing 2022-Q3 row down to row 3
#    and insert the new 2022-Q4 figures into row 2
# ---------------------------------------------------------------------
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("C3").Value = 5
$wsTotal.Range("D3").Value = 0.07000000000000001

$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 9
$wsTotal.Range("D2").Value = 0.09
